$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted for the "Feria Lagunitas de
# Puerto Montt" / Perejil series: copy the existing row 371 record and
# insert it as the new row 372, which pushes every subsequent row (the old
# 372..426) down by one (to 373..427).
$ws.Rows.Item(371).Copy()
$ws.Rows.Item(372).Insert()

# Row 371 (the original record) then gets an updated Volumen (column J)
# value, from 20 to 50.
$ws.Range("J371").Value = 50
